$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = '@'
    $c.Value = $text
    $c.Style = 'Normal'
}

Set-TextValue 'D2' '59.259.07'
Set-TextValue 'E2' '  -2.00%  '
Set-TextValue 'D3' '2.581.66'
Set-TextValue 'E3' '  -2.19%  '
Set-TextValue 'E4' '  -0.08%  '
Set-TextValue 'D5' '562.56'
Set-TextValue 'E5' '  -1.39%  '
Set-TextValue 'D6' '142.48'
Set-TextValue 'E6' '  -2.89%  '
Set-TextValue 'E7' '  +0.28%  '
Set-TextValue 'E8' '  -2.13%  '
Set-TextValue 'D9' '2.588.71'
Set-TextValue 'E9' '  -2.87%  '
Set-TextValue 'E10' '  -2.81%  '
Set-TextValue 'E11' '  -0.88%  '
Set-TextValue 'E12' '  +11.74%  '
Set-TextValue 'D13' '0.353'
Set-TextValue 'E13' '  +3.23%  '
Set-TextValue 'D14' '3.036.36'
Set-TextValue 'E14' '  -2.41%  '
Set-TextValue 'D15' '23.28'
Set-TextValue 'E15' '  +7.05%  '
Set-TextValue 'D16' '59.217.55'
Set-TextValue 'E16' '  -2.05%  '
Set-TextValue 'E17' '  -0.29%  '
Set-TextValue 'D18' '2.583.99'
Set-TextValue 'E18' '  -2.58%  '
Set-TextValue 'E19' '  +0.44%  '
Set-TextValue 'D20' '337.01'
Set-TextValue 'E20' '  -2.32%  '
Set-TextValue 'E21' '  -0.91%  '
Set-TextValue 'E22' '  -0.16%  '
Set-TextValue 'E23' '  +0.22%  '
Set-TextValue 'D24' '64.15'
Set-TextValue 'D25' '0.466'
Set-TextValue 'E25' '  +5.07%  '
Set-TextValue 'E26' '  +0.28%  '
Set-TextValue 'E27' '  -3.07%  '
Set-TextValue 'D28' '7.32'
Set-TextValue 'E28' '  -0.50%  '
Set-TextValue 'E29' '  -0.62%  '
Set-TextValue 'E30' '  +0.08%  '
Set-TextValue 'E31' '  -2.78%  '
Set-TextValue 'B32' 'Aptos'
Set-TextValue 'C32' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D32' '6.12'
Set-TextValue 'E32' '  +0.13%  '
Set-TextValue 'B33' 'Monero'
Set-TextValue 'C33' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D33' '160.41'
Set-TextValue 'E33' '  +2.77%  '
Set-TextValue 'E34' '  -1.31%  '
Set-TextValue 'D35' '4.03'
Set-TextValue 'E35' '  -1.59%  '
Set-TextValue 'E36' '  -1.06%  '
Set-TextValue 'D37' '0.877'
Set-TextValue 'E37' '  -3.82%  '
Set-TextValue 'D38' '0.871'
Set-TextValue 'E38' '  -4.40%  '
Set-TextValue 'D39' '37.43'
Set-TextValue 'E39' '  -0.49%  '
Set-TextValue 'E40' '  -2.45%  '
Set-TextValue 'B41' 'Filecoin'
Set-TextValue 'C41' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D41' '3.67'
Set-TextValue 'E41' '  -0.04%  '
Set-TextValue 'B42' 'Bittensor'
Set-TextValue 'C42' 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 'D42' '293.63'
Set-TextValue 'E42' '  -4.31%  '
Set-TextValue 'E43' '  +0.51%  '
Set-TextValue 'D44' '132.09'
Set-TextValue 'E44' '  +5.19%  '
Set-TextValue 'D45' '0.0972'
Set-TextValue 'E45' '  -0.83%  '
Set-TextValue 'D46' '0.596'
Set-TextValue 'E46' '  -2.19%  '
Set-TextValue 'E47' '  -0.09%  '
Set-TextValue 'E48' '  -2.56%  '
Set-TextValue 'E49' '  -2.43%  '
Set-TextValue 'E50' '  -1.13%  '
Set-TextValue 'E51' '  +0.12%  '
